# Update the "Förändrad" (Changed) date column (C) for every data row
# from 45209 (2023-10-10) to 45210 (2023-10-11).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)  # Column C
    if ($cell.Value2 -eq 45209) {
        $cell.Value2 = 45210
    }
}
